$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 245 (shifts old rows 245.. down by one)
$ws.Rows("245:245").Insert()

# Fill in the new row 245 with the predicted next-day entry/exit data
$newRow = 245
$ws.Cells.Item($newRow, 1).Value = 45322
$ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($newRow, 2).Value = 824.95
$ws.Cells.Item($newRow, 3).Value = 838
$ws.Cells.Item($newRow, 4).Value = 815
$ws.Cells.Item($newRow, 5).Value = 831.5
$ws.Cells.Item($newRow, 6).Value = 815.95
$ws.Cells.Item($newRow, 7).Value = 128438
$ws.Cells.Item($newRow, 8).Value = "Wednesday"
$ws.Cells.Item($newRow, 9).Value = 0.79
$ws.Cells.Item($newRow, 10).Value = 1.56
$ws.Cells.Item($newRow, 11).Value = -1.22
$ws.Cells.Item($newRow, 12).Value = 9
$ws.Cells.Item($newRow, 13).Value = -0.77
$ws.Cells.Item($newRow, 14).Value = 2.79
$ws.Cells.Item($newRow, 15).Value = 837.5715909090908
$ws.Cells.Item($newRow, 16).Value = 1
$ws.Cells.Item($newRow, 17).Value = 14
$ws.Cells.Item($newRow, 18).Value = 61
$ws.Cells.Item($newRow, 19).Value = 30
$ws.Cells.Item($newRow, 20).Value = 45
$ws.Cells.Item($newRow, 21).Value = 12
$ws.Cells.Item($newRow, 22).Value = 815
$ws.Cells.Item($newRow, 23).Value = 838
$ws.Cells.Item($newRow, 24).Value = 0
$ws.Cells.Item($newRow, 25).Value = 370
$ws.Cells.Item($newRow, 26).Value = -1.22
$ws.Cells.Item($newRow, 27).Value = 1.56
$ws.Cells.Item($newRow, 28).Value = 2.74
$ws.Cells.Item($newRow, 29).Value = 815
$ws.Cells.Item($newRow, 30).Value = 838
$ws.Cells.Item($newRow, 31).Value = 0
$ws.Cells.Item($newRow, 32).Value = 370
$ws.Cells.Item($newRow, 33).Value = -1.22
$ws.Cells.Item($newRow, 34).Value = 1.56
$ws.Cells.Item($newRow, 35).Value = 2.74

# Update the rows that shifted down (previously 245 and 246, now 246 and 247)
# Row 246 (was row 245): PvClose (F), DiffPvClose/Open (L), and 44MA (O) change
$ws.Cells.Item(246, 6).Value = 831.5
$ws.Cells.Item(246, 12).Value = -31.5
$ws.Cells.Item(246, 15).Value = 836.3681818181819

# Row 247 (was row 246): only 44MA (O) changes
$ws.Cells.Item(247, 15).Value = 835.1602272727272
